$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting the existing rows 32-34 down to 33-35.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly data entry.
$ws.Cells.Item(32, 1).Value = 9
$ws.Cells.Item(32, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(32, 3).Value = "Metropolitana"
$ws.Cells.Item(32, 4).Value = 44946
$ws.Cells.Item(32, 5).Value = 13
$ws.Cells.Item(32, 6).Value = 100112010
$ws.Cells.Item(32, 7).Value = "Achicoria"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 70
$ws.Cells.Item(32, 11).Value = 7000
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 7000
$ws.Cells.Item(32, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(32, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(32, 16).Value = 438
$ws.Cells.Item(32, 17).Value = 16
$ws.Cells.Item(32, 18).Value = "Hortaliza"
